# GSC export refresh: drop the oldest day (2025-09-30) and append the
# newest day (2025-12-29). This shifts every remaining row's date and
# "HTTPS URLs" count up by one row, and appends a fresh zeroed row for
# the new date at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the top data row (2025-09-30). This shifts rows 3..91 up to
# rows 2..90, which also moves their B/C values up by one row - exactly
# matching the diff (C2 becomes old C3, C3 becomes old C4, etc).
$ws.Rows.Item(2).Delete()

# After the delete, the table's last row (91) needs new content for the
# new date 2025-12-29, with Non-HTTPS/HTTPS URL counts starting at 0.
# A plain Value assignment of a "yyyy-MM-dd"-looking string gets
# auto-converted to a date serial by Excel, so stage the text in a
# scratch cell with a text format, copy/paste-values it into place, then
# clean the scratch cell up.
$helper = $ws.Cells.Item(200, 1)
$helper.NumberFormat = "@"
$helper.Value = "2025-12-29"
$helper.Copy()
$ws.Cells.Item(91, 1).PasteSpecial(-4163)  # xlPasteValues
$helper.Clear()

$ws.Cells.Item(91, 2).Value = 0
$ws.Cells.Item(91, 3).Value = 0
